$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center (horizontally and vertically) the used data range A1:D24
$range = $ws.Range("A1:D24")
$range.HorizontalAlignment = -4108   # xlCenter
$range.VerticalAlignment = -4108     # xlCenter

# Update selection to D2
$ws.Range("D2").Select()
